$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 63

$ws.Cells.Item($row, 1).Value = 1
$ws.Cells.Item($row, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($row, 3).Value = "Arica y Parinacota"
$ws.Cells.Item($row, 4).Value = 45239
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(62, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 15
$ws.Cells.Item($row, 6).Value = 100112044
$ws.Cells.Item($row, 7).Value = "Perejil"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 300
$ws.Cells.Item($row, 11).Value = 1300
$ws.Cells.Item($row, 12).Value = 1500
$ws.Cells.Item($row, 13).Value = 1400
$ws.Cells.Item($row, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item($row, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value = 700
$ws.Cells.Item($row, 17).Value = 2
$ws.Cells.Item($row, 18).Value = "Hortaliza"
